$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: copy the date format from B2 onto B9, then set values for B9:G9
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = 43906
$ws.Range("C9:G9").Value = 0.2

# Row 10: copy the date format from B2 onto B10, then set values for B10:G10
$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = 43913
$ws.Range("C10:G10").Value = 0.2

$excel.CutCopyMode = 0

# Update the active selection on the sheet
$ws.Range("G13").Select()
